$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 5 and 6, shifting the existing question rows
# (setID 2, 3, 4) down to rows 7-15.
$ws.Range("A5:A6").EntireRow.Insert()

# New question: setID 1, qnNo 4 (row 5)
# Fill opt4 (G5) before opt3 (F5) to reproduce the author's shared-string order.
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "In which access should a constructor be defined, so that object of the class can be created in any function?"
$ws.Range("D5").Value = "Any access specifier will work"
$ws.Range("E5").Value = "Public"
$ws.Range("G5").Value = "Private"
$ws.Range("F5").Value = "Protected"
$ws.Range("H5").Value = 2

# New question: setID 1, qnNo 5 (row 6)
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = "Which access specifier is usually used for data members of a class?"
$ws.Range("D6").Value = "Private"
$ws.Range("E6").Value = "Protected"
$ws.Range("F6").Value = "Public"
$ws.Range("G6").Value = "Default"
$ws.Range("H6").Value = 1

# Matches the selection left behind in the saved workbook.
$ws.Range("H13").Select()
